# edit.ps1 -- applies the target diff to tc_p037v.docx
#
# Summary of the change (see commit message / xml diff):
#  1) "<m>verre</m> net avecq"  ->  "<m>verre net</m> avecq"
#     (" net" moves from after </m> to before it, i.e. inside the <m> span)
#  2) "desgraisse car sil a tant soit peu de "
#        ->  "<m>desgraisse</m> car sil a tant soit peu de "
#     (wrap just the word "desgraisse" in blue Courier-New <m>/</m> markup
#      runs, matching the markup style already used elsewhere in the doc)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: move " net" from the run after </m> to the end of the "verre"
# run, immediately before </m>.
# ---------------------------------------------------------------------

$full = $d.Content.Text
$idx = $full.IndexOf("verre</m> net avecq")

if ($idx -ge 0) {
    # "verre" run occupies [$idx, $idx+5); "</m>" occupies [$idx+5, $idx+9);
    # " net avecq" occupies [$idx+9, $idx+19).

    # Remove the leading " net" (4 chars) from " net avecq", leaving " avecq".
    # Do this first since it is the later position in the document, so it
    # does not disturb the earlier offsets we still need to use.
    $netRange = $d.Range($idx + 9, $idx + 9 + 4)
    $netRange.Text = ""

    # Append " net" to the end of the "verre" run (still before </m>).
    $verreRange = $d.Range($idx, $idx + 5)
    $verreRange.InsertAfter(" net")
}

# ---------------------------------------------------------------------
# Edit 2: wrap "desgraisse" in its own <m>...</m> run pair, splitting the
# single run "desgraisse car sil a tant soit peu de " into four runs.
# ---------------------------------------------------------------------

$full = $d.Content.Text
$target = "desgraisse car sil a tant soit peu de "
$idx2 = $full.IndexOf($target)

if ($idx2 -ge 0) {
    # Use an existing "<m>" / "</m>" run elsewhere in the document as a
    # formatting template (Courier New, blue 0000ff, 9pt/sz18) so the new
    # tag runs get identical rPr to every other markup tag in the file.
    $tplOpenIdx = $full.IndexOf("<m>")
    $tplOpen = $d.Range($tplOpenIdx, $tplOpenIdx + 3)

    $tplCloseIdx = $full.IndexOf("</m>")
    $tplClose = $d.Range($tplCloseIdx, $tplCloseIdx + 4)

    # Insert the literal tag text (plain-formatted for now); insert the
    # closing tag first so the earlier offset ($idx2) stays valid.
    $afterWord = $d.Range($idx2 + 10, $idx2 + 10)   # right after "desgraisse"
    $afterWord.InsertBefore("</m>")

    $beforeWord = $d.Range($idx2, $idx2)            # right before "desgraisse"
    $beforeWord.InsertBefore("<m>")

    # Re-apply the markup formatting (font/color/size) to the two new tag
    # runs by copying the template runs' formatted text onto them.
    $openTagRange = $d.Range($idx2, $idx2 + 3)
    $openTagRange.FormattedText = $tplOpen.FormattedText

    $closeTagStart = $idx2 + 3 + 10
    $closeTagRange = $d.Range($closeTagStart, $closeTagStart + 4)
    $closeTagRange.FormattedText = $tplClose.FormattedText
}
